# Yxs_table.xlsx edit:
#  - Insert two new columns ("mu 2.1 [1/h]" and "mu 2.2 [1/h]") right after the
#    existing "mu 2 [1/h]" column (J), shifting the four "Yxs ..." columns
#    from J:M to L:O.
#  - Populate the new mu 2.1 / mu 2.2 columns: they mirror column I (mu 2),
#    split the same way the existing Yxs1/Yxs2 columns are split across the
#    data rows (mu 2.1 for rows 3-9, mu 2.2 for rows 11-21, row 10 left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at J:K - this shifts old J:M (Yxs 1.1, Yxs 1.2,
# Yxs 2.1, Yxs 2.2) to L:O automatically, carrying over styles/values.
$ws.Columns("J:K").Insert()

# New column headers
$ws.Range("J1").Value = "mu 2.1 [1/h]"
$ws.Range("K1").Value = "mu 2.2 [1/h]"

# mu 2.1 [1/h] values (copy of mu 2 for rows 3-9)
$ws.Range("J3").Value = 0.3346417012349824
$ws.Range("J4").Value = 0.8127218063500202
$ws.Range("J5").Value = 0.4243727592997257
$ws.Range("J6").Value = 0.2746536833672744
$ws.Range("J7").Value = 0.1442387591429934
$ws.Range("J8").Value = 0.2217057364726398
$ws.Range("J9").Value = 0.07735908689312881

# mu 2.2 [1/h] values (copy of mu 2 for rows 11-21; row 10 stays blank)
$ws.Range("K11").Value = 0.03719083035979676
$ws.Range("K12").Value = 0.03619552273593946
$ws.Range("K13").Value = 0.03801145000663797
$ws.Range("K14").Value = -0.06954913092679149
$ws.Range("K15").Value = -0.01970310092027705
$ws.Range("K16").Value = -0.006451635241295772
$ws.Range("K17").Value = 0.1810127751752835
$ws.Range("K18").Value = -0.02883454103336998
$ws.Range("K19").Value = 0.1444724872797511
$ws.Range("K20").Value = 0.02448929697916169
$ws.Range("K21").Value = 0.0387149826119924
